# Updates the "cryptos" price/volume snapshot to the latest scraped values.
# Note: several Price cells (column D) look like plain decimals (e.g. "522.36");
# a leading apostrophe forces Excel to keep them as literal text (matching the
# original inline-string cells) instead of auto-converting them to floating
# point numbers, which would introduce binary rounding noise.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.223.53'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.593.22'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'522.36"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Value = "'143.77"
$ws.Range("E6").Value = '  +0.85%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").Value = '2.614.67'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("D11").Value = "'0.101"
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = "'0.132"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '3.051.14'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").Value = '58.190.12'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = "'20.50"
$ws.Range("E16").Value = '  -2.10%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.610.34'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = "'338.96"
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").Value = "'10.30"
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").Value = "'6.38"
$ws.Range("E22").Value = '  +1.95%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = "'65.15"
$ws.Range("E24").Value = '  +1.56%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("E26").Value = '  -2.83%  '
$ws.Range("D27").Value = '2.723.93'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").Value = "'7.03"
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("D30").Value = '0.0₃0748'
$ws.Range("E30").Value = '  -5.17%  '
$ws.Range("D32").Value = "'6.24"
$ws.Range("E32").Value = '  -5.76%  '
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("D34").Value = "'18.80"
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").Value = "'149.83"
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("D36").Value = "'4.03"
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").Value = "'1.14"
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").Value = "'0.858"
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("D41").Value = "'36.04"
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("D42").Value = "'3.55"
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("D44").Value = "'273.13"
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Value = "'0.598"
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = "'0.0958"
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").Value = "'10.67"
$ws.Range("E47").Value = '  +0.57%  '
$ws.Range("D48").Value = "'18.79"
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").Value = "'18.97"
$ws.Range("E50").Value = '  +4.31%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = "'4.66"
$ws.Range("E51").Value = '  -0.05%  '
